$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp (17:05 -> 17:35)
$ws.Range("A1").Value = "Datos actualizados a 22 de Mayo de 2020 a las 17:35"

# Estados Unidos (row 4) - updated totals
$ws.Range("B4").Value = 1625071
$ws.Range("C4").Value = 4169
$ws.Range("E4").Value = 1145558

# Rows 19/20 - Chile overtakes Mexico in ranking.
# Chile gets new/updated numbers and moves above Mexico; Mexico keeps its
# previous numbers but drops one rank.
$ws.Range("A19").Value = "Chile"
$ws.Range("B19").Value = 61857
$ws.Range("C19").Value = 4276
$ws.Range("D19").Value = 25342
$ws.Range("E19").Value = 35885
$ws.Range("G19").Value = 41
$ws.Range("H19").Value = 630

$ws.Range("A20").Value = "Mexico"
$ws.Range("B20").Value = 59567
$ws.Range("C20").Value = 2973
$ws.Range("D20").Value = 40657
$ws.Range("E20").Value = 12400
$ws.Range("G20").Value = 420
$ws.Range("H20").Value = 6510

# Singapur (row 29) - updated active/recovered numbers
$ws.Range("D29").Value = 12995
$ws.Range("E29").Value = 17408

# Republica Dominicana (row 45) - updated totals
$ws.Range("B45").Value = 13989
$ws.Range("C45").Value = 332
$ws.Range("D45").Value = 7572
$ws.Range("E45").Value = 5961
$ws.Range("G45").Value = 8
$ws.Range("H45").Value = 456

# Argentina (row 51) - minor update
$ws.Range("E51").Value = 6480
$ws.Range("G51").Value = 3
$ws.Range("H51").Value = 419

# Rows 110-113 - Guinea Ecuatorial overtakes Mali/Niger/Republica de Chipre.
# Guinea Ecuatorial gets new/updated numbers and moves to the top of this
# group; Mali, Niger and Republica de Chipre each keep their previous
# numbers but drop one rank.
$ws.Range("A110").Value = "Guinea Ecuatorial"
$ws.Range("B110").Value = 960
$ws.Range("C110").Value = 57
$ws.Range("D110").Value = 165
$ws.Range("E110").Value = 784
$ws.Range("G110").Value = 1
$ws.Range("H110").Value = 11

$ws.Range("A111").Value = "Mali"
$ws.Range("B111").Value = 947
$ws.Range("D111").Value = 558
$ws.Range("E111").Value = 329
$ws.Range("H111").Value = 60

$ws.Range("A112").Value = "Niger"
$ws.Range("B112").Value = 924
$ws.Range("D112").Value = 753
$ws.Range("E112").Value = 111
$ws.Range("H112").Value = 60

$ws.Range("A113").Value = "Republica de Chipre"
$ws.Range("B113").Value = 923
$ws.Range("D113").Value = 561
$ws.Range("E113").Value = 345
$ws.Range("H113").Value = 17
